$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage cell used to force text-typed values via copy/paste-special (avoids Excel
# auto-converting numeric-looking strings like "510.81" into numbers).
$stage = $ws.Range("A1")

$stage.Value = "'57.841.92"
$stage.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E2").Value = "  -0.62%  "

$stage.Value = "'2.447.80"
$stage.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E3").Value = "  -1.26%  "

$ws.Range("E4").Value = "  -0.01%  "

$stage.Value = "'510.81"
$stage.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E5").Value = "  -1.97%  "

$stage.Value = "'129.71"
$stage.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E6").Value = "  -1.54%  "

$stage.Value = "'0.997"
$stage.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("E8").Value = "  -1.63%  "

$stage.Value = "'2.462.91"
$stage.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E9").Value = "  -0.85%  "

$ws.Range("E10").Value = "  -3.42%  "

$ws.Range("E11").Value = "  -0.08%  "

$stage.Value = "'5.18"
$stage.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E12").Value = "  -3.38%  "

$ws.Range("E13").Value = "  -4.58%  "

$stage.Value = "'2.882.78"
$stage.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E14").Value = "  -1.20%  "

$stage.Value = "'57.746.10"
$stage.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E15").Value = "  -0.64%  "

$stage.Value = "'21.87"
$stage.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E16").Value = "  -1.17%  "

$ws.Range("E17").Value = "  -2.25%  "

$stage.Value = "'2.459.65"
$stage.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E18").Value = "  -0.73%  "

$stage.Value = "'10.51"
$stage.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E19").Value = "  -3.18%  "

$stage.Value = "'318.10"
$stage.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E20").Value = "  -0.34%  "

$ws.Range("E21").Value = "  -1.51%  "

$ws.Range("E22").Value = "  -0.16%  "

$ws.Range("E23").Value = "  +2.98%  "

$stage.Value = "'63.00"
$stage.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E24").Value = "  -2.02%  "

$ws.Range("E25").Value = "  -2.31%  "

$stage.Value = "'0.992"
$stage.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E26").Value = "  -0.59%  "

$ws.Range("E27").Value = "  -0.28%  "

$ws.Range("E28").Value = "  -1.60%  "

$stage.Value = "'168.32"
$stage.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E29").Value = "  +0.69%  "

$ws.Range("E30").Value = "  -3.74%  "

$ws.Range("E31").Value = "  -2.55%  "

$ws.Range("E32").Value = "  -0.50%  "

$ws.Range("E33").Value = "  -3.01%  "

$ws.Range("E34").Value = "  -0.04%  "

$stage.Value = "'0.995"
$stage.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E35").Value = "  -0.24%  "

$stage.Value = "'17.73"
$stage.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E36").Value = "  -2.13%  "

$stage.Value = "'1.26"
$stage.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E37").Value = "  -3.99%  "

$stage.Value = "'3.89"
$stage.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E38").Value = "  -2.25%  "

$stage.Value = "'36.59"
$stage.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("E40").Value = "  -2.22%  "

$stage.Value = "'0.762"
$stage.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E41").Value = "  -4.00%  "

$stage.Value = "'270.68"
$stage.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E42").Value = "  -2.04%  "

$stage.Value = "'5.01"
$stage.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E43").Value = "  -0.45%  "

$stage.Value = "'3.37"
$stage.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E44").Value = "  -3.31%  "

$stage.Value = "'0.584"
$stage.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E45").Value = "  -2.06%  "

$stage.Value = "'0.0912"
$stage.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E46").Value = "  +0.67%  "

$stage.Value = "'120.22"
$stage.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E47").Value = "  -5.20%  "

$ws.Range("E48").Value = "  -0.34%  "

$stage.Value = "'17.21"
$stage.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E49").Value = "  -4.21%  "

$ws.Range("E50").Value = "  -2.29%  "

$stage.Value = "'16.66"
$stage.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$stage.Clear()
$ws.Range("E51").Value = "  -2.66%  "
